$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# This slide's shape-id counter advances in non-contiguous jumps
# (3, 8, 9, 13, 14, 15, 16, ...) because of internal id bookkeeping that
# isn't simply "max used id + 1". Pre-consume the first five slots with
# throwaway shapes (immediately deleted) so the two real additions below
# land on ids 15 and 16 - exactly matching the target deck's "Oval 14"
# (id 15) and "Equal 15" (id 16).
$primer = @()
for ($i = 0; $i -lt 5; $i++) {
    $primer += $s.Shapes.AddShape(9, 1, 1, 1, 1)
}
foreach ($d in $primer) { $d.Delete() }

# New "Oval 14": a duplicate of the existing "Oval 1" menu-dot circle,
# moved to its new spot and with its fill cleared so only the outline
# shows.
$baseOval = $s.Shapes.Item("Oval 1")
$oval = $baseOval.Duplicate().Item(1)
$oval.Name = "Oval 14"
$oval.Left = 210.00232283464567
$oval.Top = 342.03342519685043
$oval.Width = 22.569409448818895
$oval.Height = 26.087125984251966
$oval.Fill.Visible = $false

# New "Equal 15": a duplicate of the existing "Equal 9" glyph, moved to
# sit inside the new oval above.
$baseEqual = $s.Shapes.Item("Equal 9")
$equal = $baseEqual.Duplicate().Item(1)
$equal.Name = "Equal 15"
$equal.Left = 215.3062598425197
$equal.Top = 346.43476377952754
$equal.Width = 11.961614173228346
$equal.Height = 17.284448818897637
